# Auto-applies the per-cell text updates described by the commit diff
# (coin price/volume refresh + a one-row shift in the exchange-token block).
#
# All data cells on this sheet are stored as *text* (t="inlineStr"), even
# the numeric-looking Price/Volume columns (D, E). Plain COM assignment of a
# numeric-looking string (e.g. "26.50") gets auto-converted by Excel into a
# real number/percentage (losing trailing zeros / exact formatting, e.g.
# "26.50" -> 26.5, "4.700" -> 4.7000000000000002, "1.49%" -> 0.0149).
# To preserve the exact text we prefix such values with a literal leading
# apostrophe (the classic "force text" Excel entry trick) and then reset the
# cell style to "Normal" so no stray NumberFormat/quotePrefix styling is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''264.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''1.49%'
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''26.50'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''-2.18%'
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''4.700'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''0.12%'
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''0.06093'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''-1.41%'
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''6.727'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''0.67%'
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.8508'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''-0.07%'
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.9094'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''-0.32%'
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '''0.04998'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''7.22%'
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D10").Value = '''0.07107'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''0.32%'
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").Value = '''0.03146'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''1.10%'
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").Value = '''0.09028'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''-0.16%'
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").Value = '''0.001536'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''-0.41%'
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").Value = '''0.0006049'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''-2.10%'
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '''0.006049'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''-0.43%'
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '''3.449'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''-0.05%'
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '''3.169'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''0.18%'
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '''2.175'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''-0.16%'
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '''0.3072'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''-0.19%'
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = 'WazirX'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D20").Value = '''0.1409'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''0.16%'
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = '''-1.35%'
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''4.128'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''0.98%'
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''0.04247'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''0.14%'
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''0.001176'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''-3.26%'
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''0.004058'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''6.72%'
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''0.04%'
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '''6.63%'
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = '''0.03928'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''1.15%'
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.1115'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''0.45%'
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''0.004199'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''2.85%'
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''-3.54%'
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''0.01171'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''-28.30%'
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.00005097'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''-1.25%'
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''0.09%'
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = '''0.2586'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''58.17%'
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''0.09%'
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''0.09%'
$ws.Range("E50").Style = "Normal"
